$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# The edit cascades a block of text one "slot" down through several
# paragraphs/runs (each slot keeps its own paragraph/run formatting but
# receives the text that used to belong to the next slot down), finishing
# with brand-new text in the last slot, and opens up a brand-new first
# slot (a "List Bullet" paragraph under "Docente(s) Responsável(eis)")
# that receives the text which used to occupy the first slot.
#
# To keep every Find target text unique at the moment we search for it,
# we apply the cascade from the BOTTOM of the document upward (mirroring
# the diff hunks in reverse order), and only insert the brand-new
# duplicate-text bullet paragraph as the very last step.
# -----------------------------------------------------------------------

# "Norma de recuperação: " body  <-  brand-new text
$old_norma = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2"
$new_norma = "Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2"
$d.Content.Find.Execute($old_norma, $true, $true, $false, $false, $false, $true, 1, $false, $new_norma, 1) | Out-Null

# "Critério: " body  <-  old "Norma de recuperação: " body text
$old_criterio = "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T)."
$new_criterio = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2"
$d.Content.Find.Execute($old_criterio, $true, $true, $false, $false, $false, $true, 1, $false, $new_criterio, 1) | Out-Null

# "Método: " body  <-  old "Critério: " body text
$old_metodo = "1. Introdução; revisão da termodinâmica de soluções; teoria básica de equilíbrio de fases; curvas de energia livre versus composição; regra das fases; 2. Sistemas unários, equilíbrios bi-, mono- e invariantes; 3. Sistemas binários isomorfos; a regra da alavanca; solidificação em equilíbrio e fora de equilíbrio; mínimos e máximos; 4. Sistemas eutéticos binários; solidificação e microetruturas de ligas hipoeutéticas, eutéticas e hipereutéticas; solidificação unidirecional com eutéticos; casos limites de eutéticos; 5. Sistemas eutetóides binários; solidificação e microetruturas de ligas hipoeutetóides, eutetói-des e hipereutetóides; o sistema Fe-C; 6. Sistemas monotéticos; sistemas monotetóides; sistemas metatéticos; transformações congruentes; 7. Sistemas peritéticos binários; resfriamento em equilíbrio e fora do equilíbrio de ligas peritéticas; sistemas peritetóides binários; sistemas sintéticos binários; 8. Sistemas ternários isomorfos; o triângulo de Gibbs; seções isotérmicas; projeções liquidus; seções verticais; máximos e mínimos; resfriamento em equilíbrio; 9. Equilíbrio ternário de três fases; regra da alavanca em campos trifásicos; resfriamento em equilíbrio; 10. Equilíbrio ternário de quatro fases: equilíbrio de classe I; equilíbrio de classe II e equilíbrio de classe III; 11. Transformações congruentes em sistemas ternários; sistemas ternários complexos; 12. Cálculo termodinâmico de diagramas de fases."
$new_metodo = "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T)."
$d.Content.Find.Execute($old_metodo, $true, $true, $false, $false, $false, $true, 1, $false, $new_metodo, 1) | Out-Null

# "Programa" body  <-  old "Método: " body text
$old_programa = "A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários;C. Sistemas binários;D. Sistemas ternários;E. Cálculo termodinâmico de diagramas de fases."
$new_programa = "1. Introdução; revisão da termodinâmica de soluções; teoria básica de equilíbrio de fases; curvas de energia livre versus composição; regra das fases; 2. Sistemas unários, equilíbrios bi-, mono- e invariantes; 3. Sistemas binários isomorfos; a regra da alavanca; solidificação em equilíbrio e fora de equilíbrio; mínimos e máximos; 4. Sistemas eutéticos binários; solidificação e microetruturas de ligas hipoeutéticas, eutéticas e hipereutéticas; solidificação unidirecional com eutéticos; casos limites de eutéticos; 5. Sistemas eutetóides binários; solidificação e microetruturas de ligas hipoeutetóides, eutetói-des e hipereutetóides; o sistema Fe-C; 6. Sistemas monotéticos; sistemas monotetóides; sistemas metatéticos; transformações congruentes; 7. Sistemas peritéticos binários; resfriamento em equilíbrio e fora do equilíbrio de ligas peritéticas; sistemas peritetóides binários; sistemas sintéticos binários; 8. Sistemas ternários isomorfos; o triângulo de Gibbs; seções isotérmicas; projeções liquidus; seções verticais; máximos e mínimos; resfriamento em equilíbrio; 9. Equilíbrio ternário de três fases; regra da alavanca em campos trifásicos; resfriamento em equilíbrio; 10. Equilíbrio ternário de quatro fases: equilíbrio de classe I; equilíbrio de classe II e equilíbrio de classe III; 11. Transformações congruentes em sistemas ternários; sistemas ternários complexos; 12. Cálculo termodinâmico de diagramas de fases."
$d.Content.Find.Execute($old_programa, $true, $true, $false, $false, $false, $true, 1, $false, $new_programa, 1) | Out-Null

# "Programa resumido" body  <-  old "Programa" body text
$old_resumido = "5009972 - Gilberto Carvalho Coelho"
$new_resumido = "A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários;C. Sistemas binários;D. Sistemas ternários;E. Cálculo termodinâmico de diagramas de fases."
$d.Content.Find.Execute($old_resumido, $true, $true, $false, $false, $false, $true, 1, $false, $new_resumido, 1) | Out-Null

# -----------------------------------------------------------------------
# Finally, insert the brand-new "List Bullet" paragraph right after the
# "Docente(s) Responsável(eis)" heading (before "Programa resumido"),
# carrying the professor line that used to sit under "Programa resumido".
# -----------------------------------------------------------------------
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Docente(s) Responsável(eis)*") {
        $targetIndex = $i
    }
}
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter()

$newIndex = $targetIndex + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newPara.Style = "List Bullet"
$newPara.Range.Text = "5009972 - Gilberto Carvalho Coelho"
